$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.51%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.04"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.67%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.138"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.01%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07611"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.72%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.623"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.19%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2.464"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.62%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9008"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.76%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1105"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "10.36%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1771"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.63%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09230"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.65%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04194"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.45%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.48%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001251"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.72%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005844"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.75%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.07%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.266"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.86%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.94%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.544"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-7.55%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.43%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-16.34%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04128"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.95%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.36%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003999"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.55%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001301"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.31%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02396"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "2.57%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05183"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.81%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007760"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.69%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1301"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.36%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.29%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.94%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007717"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-8.95%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3058"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.81%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006733"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.32%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.06%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.01049"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "208.60%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.06%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.06%"
